$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2828.908060756878
$ws.Cells.Item(3, 2).Value = 3298.454299391392
$ws.Cells.Item(4, 2).Value = 3794.26526039282
$ws.Cells.Item(5, 2).Value = 4178.222534466338
$ws.Cells.Item(6, 2).Value = 4556.861937870579
$ws.Cells.Item(7, 2).Value = 4834.275747754319
$ws.Cells.Item(8, 2).Value = 5062.77171820571
$ws.Cells.Item(9, 2).Value = 5296.209976576868
$ws.Cells.Item(10, 2).Value = 5504.106933312641
$ws.Cells.Item(11, 2).Value = 5699.529698776136
$ws.Cells.Item(12, 2).Value = 5839.755078910467
$ws.Cells.Item(13, 2).Value = 5974.187475317275
$ws.Cells.Item(14, 2).Value = 6110.697099808131
$ws.Cells.Item(15, 2).Value = 6224.774356328168
$ws.Cells.Item(16, 2).Value = 6320.590116923513
$ws.Cells.Item(17, 2).Value = 6399.570870684403
$ws.Cells.Item(18, 2).Value = 6442.069906481417
$ws.Cells.Item(19, 2).Value = 6492.437406627102
$ws.Cells.Item(20, 2).Value = 6569.738737587023
$ws.Cells.Item(21, 2).Value = 6590.404731793462
$ws.Cells.Item(22, 2).Value = 6603.879880877508
$ws.Cells.Item(23, 2).Value = 6574.362538734874
$ws.Cells.Item(24, 2).Value = 6607.942174405302
$ws.Cells.Item(25, 2).Value = 6623.357449945452
$ws.Cells.Item(26, 2).Value = 6594.023439976461
$ws.Cells.Item(27, 2).Value = 6557.416175048669
$ws.Cells.Item(28, 2).Value = 6507.966590108575
$ws.Cells.Item(29, 2).Value = 6468.915959886022
$ws.Cells.Item(30, 2).Value = 6415.8911579546
$ws.Cells.Item(31, 2).Value = 6321.841052719808
$ws.Cells.Item(32, 2).Value = 6266.388627505204
$ws.Cells.Item(33, 2).Value = 6204.457545401357
$ws.Cells.Item(34, 2).Value = 6116.918128245487
$ws.Cells.Item(35, 2).Value = 6037.161476941037
$ws.Cells.Item(36, 2).Value = 5955.898905059021
$ws.Cells.Item(37, 2).Value = 5834.778381063357
$ws.Cells.Item(38, 2).Value = 5682.884196220623
$ws.Cells.Item(39, 2).Value = 5579.551722986069
$ws.Cells.Item(40, 2).Value = 5485.529696197539
$ws.Cells.Item(41, 2).Value = 5340.640581307429
$ws.Cells.Item(42, 2).Value = 5191.348744969853
$ws.Cells.Item(43, 2).Value = 5058.084930782884
$ws.Cells.Item(44, 2).Value = 4904.304812994588
$ws.Cells.Item(45, 2).Value = 4740.956663982107
$ws.Cells.Item(46, 2).Value = 4773.314812996503
$ws.Cells.Item(47, 2).Value = 4632.038358301177
$ws.Cells.Item(48, 2).Value = 4475.943185769258
$ws.Cells.Item(49, 2).Value = 4280.400292942837
$ws.Cells.Item(50, 2).Value = 4105.843735977068
$ws.Cells.Item(51, 2).Value = 3943.935646456168
$ws.Cells.Item(52, 2).Value = 3736.256662719521
$ws.Cells.Item(53, 2).Value = 3506.894490183995
$ws.Cells.Item(54, 2).Value = 3367.62597328793
$ws.Cells.Item(55, 2).Value = 3207.442870521121
$ws.Cells.Item(56, 2).Value = 3027.720915403213
$ws.Cells.Item(57, 2).Value = 2825.968194738774
$ws.Cells.Item(58, 2).Value = 2701.098773225646
$ws.Cells.Item(59, 2).Value = 2523.992323436909
$ws.Cells.Item(60, 2).Value = 2418.065447732264
$ws.Cells.Item(61, 2).Value = 2341.484670373488
$ws.Cells.Item(62, 2).Value = 2273.717867178393
